$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- View / pane changes ---
$ws.Application.ActiveWindow.ScrollRow = 50
$ws.Range("C50").Select()
$ws.Range("L64").Select()

# --- Data corrections (pH 5 pyrite: masses reported in mg instead of g) ---

# Row 61, 66: Include? flag flips from TRUE to FALSE
$ws.Range("M61").Value = $false
$ws.Range("M66").Value = $false

# Row 67
$ws.Range("E67").Value = 114.44991418286077
$ws.Range("F67").Value = 10.700471541638096
$ws.Range("L67").Value = 42845

# Row 68
$ws.Range("E68").Value = 735.87745446631959
$ws.Range("F68").Value = 25.04598557760059
$ws.Range("L68").Value = 42845

# Row 69
$ws.Range("E69").Value = 1324.2049550090051
$ws.Range("F69").Value = 226.27993517601215
$ws.Range("L69").Value = 42845

# Row 70
$ws.Range("E70").Value = 3793.3390740114323
$ws.Range("F70").Value = 436.9919422995963
$ws.Range("L70").Value = 42845

# Row 71
$ws.Range("E71").Value = 7186.3934914274823
$ws.Range("F71").Value = 1426.1909061889894
$ws.Range("L71").Value = 42845
$ws.Range("M71").Value = $false

# Row 76, 81: Include? flag flips from TRUE to FALSE
$ws.Range("M76").Value = $false
$ws.Range("M81").Value = $false
